$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 19:04"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1299772
$ws.Range("C4").Value = 7149
$ws.Range("D4").Value = 219485
$ws.Range("E4").Value = 1002728

# Row 11 - Brasil
$ws.Range("B11").Value = 140023
$ws.Range("C11").Value = 4330
$ws.Range("E11").Value = 75073
$ws.Range("G11").Value = 412
$ws.Range("H11").Value = 9600

# Row 12 - Turquia
$ws.Range("B12").Value = 135569
$ws.Range("C12").Value = 1848
$ws.Range("D12").Value = 86396
$ws.Range("E12").Value = 45484
$ws.Range("F12").Value = 1219
$ws.Range("G12").Value = 48
$ws.Range("H12").Value = 3689

# Rows 192/193 - swap Nueva Caledonia and Belice (both label + stats)
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
